$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first")

# Row 2
$ws.Range("B2").Value = 0.1557069725526251
$ws.Range("C2").Value = 0.9141480480978561
$ws.Range("D2").Value = 1.710317365282914
$ws.Range("E2").Value = 1.307791025081191
$ws.Range("F2").Value = 1.313161276524253
$ws.Range("G2").Value = 45

# Row 3
$ws.Range("B3").Value = 0.2469628145469965
$ws.Range("C3").Value = 1.238503523069014
$ws.Range("D3").Value = 3.975336168247188
$ws.Range("E3").Value = 1.993824507886085
$ws.Range("F3").Value = 1.985678071014529

# Row 4
$ws.Range("B4").Value = 0.2987143199435805
$ws.Range("C4").Value = 1.355676670930159
$ws.Range("D4").Value = 8.587743180333536
$ws.Range("E4").Value = 2.930485144192602
$ws.Range("F4").Value = 2.937560091001783
$ws.Range("G4").Value = 66
